# Apply "Improve params default settings" edit to the OpenBB Portfolio Template workbook.

$wb = $excel.ActiveWorkbook

$wsOpt = $wb.Worksheets.Item("Optimization")
$wsCfg = $wb.Worksheets.Item("CONFIG")

# 1. Unhide the CONFIG sheet.
$wsCfg.Visible = -1

# 2. Change the selected "technique" parameter from "maxsharpe" to "blacklitterman".
$wsOpt.Range("C16").Value = "blacklitterman"

# 3. Update the CONFIG table so that the "objective" parameter is marked as
#    applicable ("YES") for the "blacklitterman" technique column (M32).
$wsCfg.Range("M32").Value = "YES"

# 4. Restore the active selections on each sheet.
[void]$wsOpt.Range("C55").Select()
[void]$wsCfg.Range("M33").Select()

# Make sure the Optimization sheet is the active / displayed sheet.
[void]$wsOpt.Activate()
